# 5/23 Working time update
$wb = $excel.ActiveWorkbook

$wsPipette = $wb.Worksheets.Item("Pipette")
$wsLFGen   = $wb.Worksheets.Item("LF Gen")

# ---------------------------------------------------------------------------
# LF Gen sheet: row 14 (item #9) gets a new "Gen Main / Old" request entry.
# ---------------------------------------------------------------------------
$wsLFGen.Range("C14").Value = "Gen Main"
$wsLFGen.Range("D14").Value = "Old"
$wsLFGen.Range("E14").Value = 43243
$wsLFGen.Range("G14").Value = "OPEN"
$wsLFGen.Range("H14").Value = "Gen Main 부품 요청 - 5ea 기준`nGen Main bare PCB 요청 - 2ea"

# Give H14 the wrap-text format already used for similar two-line notes
# (same style as Pipette!H12), then bump the row height to fit it.
$wsPipette.Range("H12").Copy() | Out-Null
$wsLFGen.Range("H14").PasteSpecial(-4122) | Out-Null
$wsLFGen.Range("H14").Value = "Gen Main 부품 요청 - 5ea 기준`nGen Main bare PCB 요청 - 2ea"
$wsLFGen.Rows.Item(14).RowHeight = 33

# ---------------------------------------------------------------------------
# Pipette sheet: row 13 (item #8) is now Closed, and row 14 (item #9) gets a
# new Transformer / V2.0 request entry.
# ---------------------------------------------------------------------------

# Row 13 becomes a shaded "closed" row - copy that look from an existing
# closed-row entry (LF Gen row 9) onto Pipette B13:I13, then fix the values.
$wsLFGen.Range("B9:I9").Copy() | Out-Null
$wsPipette.Range("B13:I13").PasteSpecial(-4122) | Out-Null

$wsPipette.Range("G13").Value = "Closed"
$wsPipette.Range("I13").Value = "[5/17 Femto] - No molding 8ea 입고"

# Row 14: new Transformer / V2.0 / OPEN entry, issued 5/23.
$wsPipette.Range("C14").Value = "Transformer"
$wsPipette.Range("D14").Value = "V2.0"
$wsPipette.Range("E14").Value = 43243
$wsPipette.Range("G14").Value = "OPEN"
$wsPipette.Range("H14").Value = "Molding된 Transformer 요청 - 각 2ea씩"

# ---------------------------------------------------------------------------
# Sheet views: LF Gen becomes the active/selected tab; selections move too.
# ---------------------------------------------------------------------------
$wsPipette.Range("F19").Select() | Out-Null
$wsLFGen.Activate() | Out-Null
$wsLFGen.Range("H20").Select() | Out-Null

# ---------------------------------------------------------------------------
# Workbook calculation mode -> manual.
# ---------------------------------------------------------------------------
$excel.Calculation = -4135
